# Natmi following Dr Hou advice
# Updates the Icosl-Cd28 LR-pair sheet: rows 2-11 get new computed
# statistics (ligand/receptor detection, specificity and edge weights),
# and rows 12-16 are appended as new sender/target cluster combinations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 ---
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 0.7779683333333334
$ws.Cells.Item(2, 8).Value = 2.333905
$ws.Cells.Item(2, 9).Value = 0.03257389909260204
$ws.Cells.Item(2, 10).Value = 0.03376044640127995
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 3.473964666666667
$ws.Cells.Item(2, 14).Value = 10.421894
$ws.Cells.Item(2, 15).Value = 0.5449432418219904
$ws.Cells.Item(2, 16).Value = 0.5462679851731443
$ws.Cells.Item(2, 17).Value = 2.702634501785556
$ws.Cells.Item(2, 18).Value = 24.32371051607
$ws.Cells.Item(2, 19).Value = 0.01775092617030495
$ws.Cells.Item(2, 20).Value = 0.01844225103417313
# --- row 3 ---
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 0.7779683333333334
$ws.Cells.Item(3, 8).Value = 2.333905
$ws.Cells.Item(3, 9).Value = 0.03257389909260204
$ws.Cells.Item(3, 10).Value = 0.03376044640127995
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 2.854567333333333
$ws.Cells.Item(3, 14).Value = 8.563702
$ws.Cells.Item(3, 15).Value = 0.4477815193550676
$ws.Cells.Item(3, 16).Value = 0.4488700649961731
$ws.Cells.Item(3, 17).Value = 2.220762990701111
$ws.Cells.Item(3, 18).Value = 19.98686691631
$ws.Cells.Item(3, 19).Value = 0.014585990027004
$ws.Cells.Item(3, 20).Value = 0.01515405377044235
# --- row 4 ---
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 0.7779683333333334
$ws.Cells.Item(4, 8).Value = 2.333905
$ws.Cells.Item(4, 9).Value = 0.03257389909260204
$ws.Cells.Item(4, 10).Value = 0.03376044640127995
$ws.Cells.Item(4, 12).Value = 0.5
$ws.Cells.Item(4, 13).Value = 0.046379
$ws.Cells.Item(4, 14).Value = 0.092758
$ws.Cells.Item(4, 15).Value = 0.007275238822941998
$ws.Cells.Item(4, 16).Value = 0.004861949830682458
$ws.Cells.Item(4, 17).Value = 0.03608139333166666
$ws.Cells.Item(4, 18).Value = 0.21648835999
$ws.Cells.Item(4, 19).Value = 0.0002369828952930935
$ws.Cells.Item(4, 20).Value = 0.0001641415966644673
# --- row 5 ---
$ws.Cells.Item(5, 4).Value = "M1"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 0.9308339999999999
$ws.Cells.Item(5, 8).Value = 2.792502
$ws.Cells.Item(5, 9).Value = 0.03897445627130897
$ws.Cells.Item(5, 10).Value = 0.04039415233116475
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 3.473964666666667
$ws.Cells.Item(5, 14).Value = 10.421894
$ws.Cells.Item(5, 15).Value = 0.5449432418219904
$ws.Cells.Item(5, 16).Value = 0.5462679851731443
$ws.Cells.Item(5, 17).Value = 3.233684426532
$ws.Cells.Item(5, 18).Value = 29.103159838788
$ws.Cells.Item(5, 19).Value = 0.02123886654873652
$ws.Cells.Item(5, 20).Value = 0.02206603220672244
# --- row 6 ---
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 0.9308339999999999
$ws.Cells.Item(6, 8).Value = 2.792502
$ws.Cells.Item(6, 9).Value = 0.03897445627130897
$ws.Cells.Item(6, 10).Value = 0.04039415233116475
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 2.854567333333333
$ws.Cells.Item(6, 14).Value = 8.563702
$ws.Cells.Item(6, 15).Value = 0.4477815193550676
$ws.Cells.Item(6, 16).Value = 0.4488700649961731
$ws.Cells.Item(6, 17).Value = 2.657128329156
$ws.Cells.Item(6, 18).Value = 23.91415496240399
$ws.Cells.Item(6, 19).Value = 0.01745204124520438
$ws.Cells.Item(6, 20).Value = 0.01813172578235524
# --- row 7 ---
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 0.9308339999999999
$ws.Cells.Item(7, 8).Value = 2.792502
$ws.Cells.Item(7, 9).Value = 0.03897445627130897
$ws.Cells.Item(7, 10).Value = 0.04039415233116475
$ws.Cells.Item(7, 12).Value = 0.5
$ws.Cells.Item(7, 13).Value = 0.046379
$ws.Cells.Item(7, 14).Value = 0.092758
$ws.Cells.Item(7, 15).Value = 0.007275238822941998
$ws.Cells.Item(7, 16).Value = 0.004861949830682458
$ws.Cells.Item(7, 17).Value = 0.043171150086
$ws.Cells.Item(7, 18).Value = 0.259026900516
$ws.Cells.Item(7, 19).Value = 0.0002835484773680823
$ws.Cells.Item(7, 20).Value = 0.0001963943420870679
# --- row 8 ---
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 11.921496
$ws.Cells.Item(8, 8).Value = 35.764488
$ws.Cells.Item(8, 9).Value = 0.4991586303686639
$ws.Cells.Item(8, 10).Value = 0.5173411429313618
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 3.473964666666667
$ws.Cells.Item(8, 14).Value = 10.421894
$ws.Cells.Item(8, 15).Value = 0.5449432418219904
$ws.Cells.Item(8, 16).Value = 0.5462679851731443
$ws.Cells.Item(8, 17).Value = 41.414855877808
$ws.Cells.Item(8, 18).Value = 372.733702900272
$ws.Cells.Item(8, 19).Value = 0.2720131222165243
$ws.Cells.Item(8, 20).Value = 0.2826069037962867
# --- row 9 ---
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 11.921496
$ws.Cells.Item(9, 8).Value = 35.764488
$ws.Cells.Item(9, 9).Value = 0.4991586303686639
$ws.Cells.Item(9, 10).Value = 0.5173411429313618
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 2.854567333333333
$ws.Cells.Item(9, 14).Value = 8.563702
$ws.Cells.Item(9, 15).Value = 0.4477815193550676
$ws.Cells.Item(9, 16).Value = 0.4488700649961731
$ws.Cells.Item(9, 17).Value = 34.03071304606399
$ws.Cells.Item(9, 18).Value = 306.276417414576
$ws.Cells.Item(9, 19).Value = 0.2235140099056749
$ws.Cells.Item(9, 20).Value = 0.2322189524527948
# --- row 10 ---
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 11.921496
$ws.Cells.Item(10, 8).Value = 35.764488
$ws.Cells.Item(10, 9).Value = 0.4991586303686639
$ws.Cells.Item(10, 10).Value = 0.5173411429313618
$ws.Cells.Item(10, 12).Value = 0.5
$ws.Cells.Item(10, 13).Value = 0.046379
$ws.Cells.Item(10, 14).Value = 0.092758
$ws.Cells.Item(10, 15).Value = 0.007275238822941998
$ws.Cells.Item(10, 16).Value = 0.004861949830682458
$ws.Cells.Item(10, 17).Value = 0.5529070629839999
$ws.Cells.Item(10, 18).Value = 3.317442377904
$ws.Cells.Item(10, 19).Value = 0.003631498246464658
$ws.Cells.Item(10, 20).Value = 0.002515286682280204
# --- row 11 ---
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 4).Value = "M1"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 7.734682333333333
$ws.Cells.Item(11, 8).Value = 23.204047
$ws.Cells.Item(11, 9).Value = 0.3238547779442587
$ws.Cells.Item(11, 10).Value = 0.3356516160838941
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 3.473964666666667
$ws.Cells.Item(11, 14).Value = 10.421894
$ws.Cells.Item(11, 15).Value = 0.5449432418219904
$ws.Cells.Item(11, 16).Value = 0.5462679851731443
$ws.Cells.Item(11, 17).Value = 26.87001313389089
$ws.Cells.Item(11, 18).Value = 241.830118205018
$ws.Cells.Item(11, 19).Value = 0.1764824725724851
$ws.Cells.Item(11, 20).Value = 0.1833557320382586
# --- row 12 ---
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Icosl"
$ws.Cells.Item(12, 3).Value = "Cd28"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 7.734682333333333
$ws.Cells.Item(12, 8).Value = 23.204047
$ws.Cells.Item(12, 9).Value = 0.3238547779442587
$ws.Cells.Item(12, 10).Value = 0.3356516160838941
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.854567333333333
$ws.Cells.Item(12, 14).Value = 8.563702
$ws.Cells.Item(12, 15).Value = 0.4477815193550676
$ws.Cells.Item(12, 16).Value = 0.4488700649961731
$ws.Cells.Item(12, 17).Value = 22.07917152244378
$ws.Cells.Item(12, 18).Value = 198.712543701994
$ws.Cells.Item(12, 19).Value = 0.1450161845182782
$ws.Cells.Item(12, 20).Value = 0.1506639627276481
# --- row 13 ---
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Icosl"
$ws.Cells.Item(13, 3).Value = "Cd28"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 7.734682333333333
$ws.Cells.Item(13, 8).Value = 23.204047
$ws.Cells.Item(13, 9).Value = 0.3238547779442587
$ws.Cells.Item(13, 10).Value = 0.3356516160838941
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.5
$ws.Cells.Item(13, 13).Value = 0.046379
$ws.Cells.Item(13, 14).Value = 0.092758
$ws.Cells.Item(13, 15).Value = 0.007275238822941998
$ws.Cells.Item(13, 16).Value = 0.004861949830682458
$ws.Cells.Item(13, 17).Value = 0.3587268319376666
$ws.Cells.Item(13, 18).Value = 2.152360991626
$ws.Cells.Item(13, 19).Value = 0.002356120853495331
$ws.Cells.Item(13, 20).Value = 0.001631921317987382
# --- row 14 ---
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Icosl"
$ws.Cells.Item(14, 3).Value = "Cd28"
$ws.Cells.Item(14, 4).Value = "M1"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 2.5182005
$ws.Cells.Item(14, 8).Value = 5.036401
$ws.Cells.Item(14, 9).Value = 0.1054382363231665
$ws.Cells.Item(14, 10).Value = 0.07285264225229936
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.473964666666667
$ws.Cells.Item(14, 14).Value = 10.421894
$ws.Cells.Item(14, 15).Value = 0.5449432418219904
$ws.Cells.Item(14, 16).Value = 0.5462679851731443
$ws.Cells.Item(14, 17).Value = 8.748139560582333
$ws.Cells.Item(14, 18).Value = 52.488837363494
$ws.Cells.Item(14, 19).Value = 0.05745785431393952
$ws.Cells.Item(14, 20).Value = 0.03979706609770346
# --- row 15 ---
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Icosl"
$ws.Cells.Item(15, 3).Value = "Cd28"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 2.5182005
$ws.Cells.Item(15, 8).Value = 5.036401
$ws.Cells.Item(15, 9).Value = 0.1054382363231665
$ws.Cells.Item(15, 10).Value = 0.07285264225229936
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.854567333333333
$ws.Cells.Item(15, 14).Value = 8.563702
$ws.Cells.Item(15, 15).Value = 0.4477815193550676
$ws.Cells.Item(15, 16).Value = 0.4488700649961731
$ws.Cells.Item(15, 17).Value = 7.188372886083666
$ws.Cells.Item(15, 18).Value = 43.13023731650199
$ws.Cells.Item(15, 19).Value = 0.04721329365890619
$ws.Cells.Item(15, 20).Value = 0.03270137026293256
# --- row 16 ---
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Icosl"
$ws.Cells.Item(16, 3).Value = "Cd28"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 2.5182005
$ws.Cells.Item(16, 8).Value = 5.036401
$ws.Cells.Item(16, 9).Value = 0.1054382363231665
$ws.Cells.Item(16, 10).Value = 0.07285264225229936
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.5
$ws.Cells.Item(16, 13).Value = 0.046379
$ws.Cells.Item(16, 14).Value = 0.092758
$ws.Cells.Item(16, 15).Value = 0.007275238822941998
$ws.Cells.Item(16, 16).Value = 0.004861949830682458
$ws.Cells.Item(16, 17).Value = 0.1167916209895
$ws.Cells.Item(16, 18).Value = 0.4671664839579999
$ws.Cells.Item(16, 19).Value = 0.0007670883503208344
$ws.Cells.Item(16, 20).Value = 0.0003542058916633365
